$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'41.268.47"
$ws.Range("E2").Value = "  -3.85%  "
$ws.Range("D3").Value = "'2.466.72"
$ws.Range("E3").Value = "  -2.70%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "'311.44"
$ws.Range("E5").Value = "  +0.13%  "
$ws.Range("D6").Value = "'94.61"
$ws.Range("E6").Value = "  -6.06%  "
$ws.Range("D7").Value = "'0.548"
$ws.Range("E7").Value = "  -3.39%  "
$ws.Range("E8").Value = "  -0.12%  "
$ws.Range("E9").Value = "  -4.53%  "
$ws.Range("D10").Value = "'33.52"
$ws.Range("E10").Value = "  -6.41%  "
$ws.Range("E11").Value = "  -3.12%  "
$ws.Range("E12").Value = "  -1.14%  "
$ws.Range("E13").Value = "  -4.38%  "
$ws.Range("D14").Value = "'2.846.41"
$ws.Range("E14").Value = "  -2.71%  "
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").Value = "'14.99"
$ws.Range("E15").Value = "  -3.20%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "'2.440.40"
$ws.Range("E16").Value = "  -4.26%  "
$ws.Range("D17").Value = "'0.787"
$ws.Range("E17").Value = "  -3.83%  "
$ws.Range("D18").Value = "'41.291.19"
$ws.Range("E18").Value = "  -3.71%  "
$ws.Range("D19").Value = "'6.32"
$ws.Range("E19").Value = "  -5.52%  "
$ws.Range("E20").Value = "  -3.09%  "
$ws.Range("D21").Value = "'11.20"
$ws.Range("E21").Value = "  -9.65%  "
$ws.Range("D22").Value = "'68.59"
$ws.Range("E22").Value = "  -1.79%  "
$ws.Range("D23").Value = "'236.88"
$ws.Range("E23").Value = "  -2.92%  "
$ws.Range("D24").Value = "'2.75"
$ws.Range("E24").Value = "  -4.80%  "
$ws.Range("E25").Value = "  -0.04%  "
$ws.Range("E26").Value = "  -7.12%  "
$ws.Range("D27").Value = "'24.08"
$ws.Range("E27").Value = "  -6.07%  "
$ws.Range("D28").Value = "'2.23"
$ws.Range("E28").Value = "  -4.68%  "
$ws.Range("D29").Value = "'9.66"
$ws.Range("E29").Value = "  -5.61%  "
$ws.Range("D30").Value = "'36.51"
$ws.Range("E30").Value = "  -6.18%  "
$ws.Range("D31").Value = "'151.76"
$ws.Range("E31").Value = "  -5.34%  "
$ws.Range("D32").Value = "'5.49"
$ws.Range("E32").Value = "  -5.81%  "
$ws.Range("E33").Value = "  -4.42%  "
$ws.Range("D34").Value = "'2.57"
$ws.Range("E34").Value = "  -3.75%  "
$ws.Range("D35").Value = "'0.0750"
$ws.Range("E35").Value = "  -5.47%  "
$ws.Range("E36").Value = "  -3.49%  "
$ws.Range("D37").Value = "'17.16"
$ws.Range("E37").Value = "  -6.49%  "
$ws.Range("D38").Value = "'1.87"
$ws.Range("E38").Value = "  -5.43%  "
$ws.Range("E39").Value = "  -3.13%  "
$ws.Range("E40").Value = "  -7.88%  "
$ws.Range("D41").Value = "'4.22"
$ws.Range("E41").Value = "  +0.71%  "
$ws.Range("E42").Value = "  -0.04%  "
$ws.Range("D43").Value = "'19.75"
$ws.Range("E43").Value = "  -9.47%  "
$ws.Range("D44").Value = "'1.983.51"
$ws.Range("E44").Value = "  -1.03%  "
$ws.Range("E45").Value = "  -4.76%  "
$ws.Range("E46").Value = "  -9.25%  "
$ws.Range("D47").Value = "'8.71"
$ws.Range("E47").Value = "  -6.34%  "
$ws.Range("D48").Value = "'2.711.76"
$ws.Range("E48").Value = "  -2.36%  "
$ws.Range("D49").Value = "'69.61"
$ws.Range("E49").Value = "  -4.12%  "
$ws.Range("D50").Value = "'96.47"
$ws.Range("E50").Value = "  -5.03%  "
$ws.Range("B51").Value = "BitcoinSV"
$ws.Range("C51").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D51").Value = "'74.66"
$ws.Range("E51").Value = "  -6.53%  "
